# Generate Report for Handoff
# Adds a new tracked file (f4a69fa1-9443-4784-8c9a-8a524f3fdc3f.md) with status
# "Handoff transform failed" right before the existing ff090798-...md row, on
# all three worksheets (Overview, zh-cn, de-de). This pushes the
# ff090798-...md and .localization-config rows down by one.

$wb = $excel.ActiveWorkbook

$HyperlinkBlue = 15570276   # RGB(100,149,237) == 0x6495ED, packed as BGR long for COM

function Fix-HyperlinkStyle($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $HyperlinkBlue
}

# ---------------------------------------------------------------------------
# Sheet 1: Overview  (columns: A=File Name, B=zh-cn, C=de-de)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Make room: push current row 7 (ff090798-...) and row 8 (.localization-config)
# down by one row.
$ws1.Rows.Item(7).Insert()

# New row 7: the newly handed-off file.
$ws1.Range("A7").Value = "f4a69fa1-9443-4784-8c9a-8a524f3fdc3f.md"
$ws1.Range("B7").Value = "Handoff transform failed"
$ws1.Range("C7").Value = "Handoff transform failed"

# Refresh dimension-affecting content on the (shifted) rows 8 and 9 so the
# values match what was already there (Insert already moved them, this is a
# no-op safety net in case the engine ever changes behavior).
$ws1.Range("A8").Value = "ff090798-82ce-4771-adaf-679755eac184.md"
$ws1.Range("B8").Value = "Handoff transform failed"
$ws1.Range("C8").Value = "Handoff transform failed"

$ws1.Range("A9").Value = ".localization-config"
$ws1.Range("B9").Value = "Not to be localized"
$ws1.Range("C9").Value = "Not to be localized"

# Hyperlinks are not auto-shifted by Insert() in this engine, so rebuild them.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/9f8cb374d761410977483bf91ae08fca3f7a33b4/e2e/79394527-517c-49e5-98c7-f56ee5a2bdfe.md", "", "", "79394527-517c-49e5-98c7-f56ee5a2bdfe.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c8b22d56d51f62d40cfe916715c117eaca010c91/e2e/7bae3312-cc79-4c8e-a8a9-1d76a4cd91fd.md", "", "", "7bae3312-cc79-4c8e-a8a9-1d76a4cd91fd.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f42f561feddbc8605eff51a630ce818c1b691d6c/e2e/88ec1a48-be7e-4138-afc3-e1d7a68b4fa7.md", "", "", "88ec1a48-be7e-4138-afc3-e1d7a68b4fa7.md")
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/c8b22d56d51f62d40cfe916715c117eaca010c91/e2e/8a407dfc-124d-44cb-957a-d307955c1e31.md", "", "", "8a407dfc-124d-44cb-957a-d307955c1e31.md")
$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/b92849358834d0f3fa7abe8a80dea8e2a79b455a/e2e/b45e65bf-9b00-4f3c-95a2-156c21f43fc7.md", "", "", "b45e65bf-9b00-4f3c-95a2-156c21f43fc7.md")
$ws1.Hyperlinks.Add($ws1.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/9f8cb374d761410977483bf91ae08fca3f7a33b4/e2e/f4a69fa1-9443-4784-8c9a-8a524f3fdc3f.md", "", "", "f4a69fa1-9443-4784-8c9a-8a524f3fdc3f.md")
$ws1.Hyperlinks.Add($ws1.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/f42f561feddbc8605eff51a630ce818c1b691d6c/e2e/ff090798-82ce-4771-adaf-679755eac184.md", "", "", "ff090798-82ce-4771-adaf-679755eac184.md")
$ws1.Hyperlinks.Add($ws1.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/9f8cb374d761410977483bf91ae08fca3f7a33b4/.localization-config", "", "", ".localization-config")

for ($r = 2; $r -le 9; $r++) {
    Fix-HyperlinkStyle $ws1.Range("A$r")
}

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn  (A=Source File Name, B=Status, D=Latest Handoff Datetime,
#                  G=Latest Handback DateTime, H=Handoff Reason)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Rows.Item(7).Insert()

$ws2.Range("A7").Value = "f4a69fa1-9443-4784-8c9a-8a524f3fdc3f.md"
$ws2.Range("B7").Value = "Handoff transform failed"
$ws2.Range("D7").Value = "0001-01-01 00:00:00"
$ws2.Range("D7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("G7").Value = "0001-01-01 00:00:00"
$ws2.Range("H7").Value = "Ignored"

$ws2.Range("A8").Value = "ff090798-82ce-4771-adaf-679755eac184.md"
$ws2.Range("B8").Value = "Handoff transform failed"
$ws2.Range("D8").Value = "0001-01-01 00:00:00"
$ws2.Range("G8").Value = "0001-01-01 00:00:00"
$ws2.Range("H8").Value = "Ignored"

$ws2.Range("A9").Value = ".localization-config"
$ws2.Range("B9").Value = "Not to be localized"
$ws2.Range("D9").Value = "0001-01-01 00:00:00"
$ws2.Range("G9").Value = "0001-01-01 00:00:00"
$ws2.Range("H9").Value = "Ignored"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/9f8cb374d761410977483bf91ae08fca3f7a33b4/e2e/79394527-517c-49e5-98c7-f56ee5a2bdfe.md", "", "", "79394527-517c-49e5-98c7-f56ee5a2bdfe.md")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c8b22d56d51f62d40cfe916715c117eaca010c91/e2e/7bae3312-cc79-4c8e-a8a9-1d76a4cd91fd.md", "", "", "7bae3312-cc79-4c8e-a8a9-1d76a4cd91fd.md")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f42f561feddbc8605eff51a630ce818c1b691d6c/e2e/88ec1a48-be7e-4138-afc3-e1d7a68b4fa7.md", "", "", "88ec1a48-be7e-4138-afc3-e1d7a68b4fa7.md")
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/c8b22d56d51f62d40cfe916715c117eaca010c91/e2e/8a407dfc-124d-44cb-957a-d307955c1e31.md", "", "", "8a407dfc-124d-44cb-957a-d307955c1e31.md")
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/b92849358834d0f3fa7abe8a80dea8e2a79b455a/e2e/b45e65bf-9b00-4f3c-95a2-156c21f43fc7.md", "", "", "b45e65bf-9b00-4f3c-95a2-156c21f43fc7.md")
$ws2.Hyperlinks.Add($ws2.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/9f8cb374d761410977483bf91ae08fca3f7a33b4/e2e/f4a69fa1-9443-4784-8c9a-8a524f3fdc3f.md", "", "", "f4a69fa1-9443-4784-8c9a-8a524f3fdc3f.md")
$ws2.Hyperlinks.Add($ws2.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/f42f561feddbc8605eff51a630ce818c1b691d6c/e2e/ff090798-82ce-4771-adaf-679755eac184.md", "", "", "ff090798-82ce-4771-adaf-679755eac184.md")
$ws2.Hyperlinks.Add($ws2.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/9f8cb374d761410977483bf91ae08fca3f7a33b4/.localization-config", "", "", ".localization-config")

for ($r = 2; $r -le 9; $r++) {
    Fix-HyperlinkStyle $ws2.Range("A$r")
}

# ---------------------------------------------------------------------------
# Sheet 3: de-de  (same layout as zh-cn)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Rows.Item(7).Insert()

$ws3.Range("A7").Value = "f4a69fa1-9443-4784-8c9a-8a524f3fdc3f.md"
$ws3.Range("B7").Value = "Handoff transform failed"
$ws3.Range("D7").Value = "0001-01-01 00:00:00"
$ws3.Range("D7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("G7").Value = "0001-01-01 00:00:00"
$ws3.Range("H7").Value = "Ignored"

$ws3.Range("A8").Value = "ff090798-82ce-4771-adaf-679755eac184.md"
$ws3.Range("B8").Value = "Handoff transform failed"
$ws3.Range("D8").Value = "0001-01-01 00:00:00"
$ws3.Range("G8").Value = "0001-01-01 00:00:00"
$ws3.Range("H8").Value = "Ignored"

$ws3.Range("A9").Value = ".localization-config"
$ws3.Range("B9").Value = "Not to be localized"
$ws3.Range("D9").Value = "0001-01-01 00:00:00"
$ws3.Range("G9").Value = "0001-01-01 00:00:00"
$ws3.Range("H9").Value = "Ignored"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/9f8cb374d761410977483bf91ae08fca3f7a33b4/e2e/79394527-517c-49e5-98c7-f56ee5a2bdfe.md", "", "", "79394527-517c-49e5-98c7-f56ee5a2bdfe.md")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c8b22d56d51f62d40cfe916715c117eaca010c91/e2e/7bae3312-cc79-4c8e-a8a9-1d76a4cd91fd.md", "", "", "7bae3312-cc79-4c8e-a8a9-1d76a4cd91fd.md")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f42f561feddbc8605eff51a630ce818c1b691d6c/e2e/88ec1a48-be7e-4138-afc3-e1d7a68b4fa7.md", "", "", "88ec1a48-be7e-4138-afc3-e1d7a68b4fa7.md")
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/c8b22d56d51f62d40cfe916715c117eaca010c91/e2e/8a407dfc-124d-44cb-957a-d307955c1e31.md", "", "", "8a407dfc-124d-44cb-957a-d307955c1e31.md")
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/b92849358834d0f3fa7abe8a80dea8e2a79b455a/e2e/b45e65bf-9b00-4f3c-95a2-156c21f43fc7.md", "", "", "b45e65bf-9b00-4f3c-95a2-156c21f43fc7.md")
$ws3.Hyperlinks.Add($ws3.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/9f8cb374d761410977483bf91ae08fca3f7a33b4/e2e/f4a69fa1-9443-4784-8c9a-8a524f3fdc3f.md", "", "", "f4a69fa1-9443-4784-8c9a-8a524f3fdc3f.md")
$ws3.Hyperlinks.Add($ws3.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/f42f561feddbc8605eff51a630ce818c1b691d6c/e2e/ff090798-82ce-4771-adaf-679755eac184.md", "", "", "ff090798-82ce-4771-adaf-679755eac184.md")
$ws3.Hyperlinks.Add($ws3.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/9f8cb374d761410977483bf91ae08fca3f7a33b4/.localization-config", "", "", ".localization-config")

for ($r = 2; $r -le 9; $r++) {
    Fix-HyperlinkStyle $ws3.Range("A$r")
}

Write-Host "Done: inserted f4a69fa1-9443-4784-8c9a-8a524f3fdc3f.md row on all sheets."
